$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.003254057381280262
$arr[0,1] = 0.03349782192518092
$arr[0,2] = 0.4333214853875944
$arr[0,3] = 1.479172793580574
$arr[0,4] = 1.560190433612945
$arr[0,5] = 1.022540459300615
$arr[0,6] = 1.136217968661228
$ws.Range("C2:I2").Value = $arr
$ws.Range("N2").Value = 2.873581687696344

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.002821308720957205
$arr[0,1] = 0.02934618752100704
$arr[0,2] = 0.3770546871429588
$arr[0,3] = 1.338957110356432
$arr[0,4] = 1.390408734690794
$arr[0,5] = 0.94984346182747
$arr[0,6] = 1.027328598235442
$ws.Range("C3:I3").Value = $arr
$ws.Range("N3").Value = 2.562605684679454

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.002556772847142952
$arr[0,1] = 0.02681457152093003
$arr[0,2] = 0.3427093422305489
$arr[0,3] = 1.253856549063869
$arr[0,4] = 1.287088664196489
$arr[0,5] = 0.9059721192090819
$arr[0,6] = 0.9612294625744369
$ws.Range("C4:I4").Value = $arr
$ws.Range("N4").Value = 2.371325805375761

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.002449222315316746
$arr[0,1] = 0.02578698926833312
$arr[0,2] = 0.3287593870395966
$arr[0,3] = 1.219418417453937
$arr[0,4] = 1.245207559488222
$arr[0,5] = 0.8882818619925104
$arr[0,6] = 0.9344778148814328
$ws.Range("C5:I5").Value = $arr
$ws.Range("N5").Value = 2.293303068605894

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.002431377428539605
$arr[0,1] = 0.02561659642684333
$arr[0,2] = 0.32644564704799
$arr[0,3] = 1.213714304211692
$arr[0,4] = 1.238266383928618
$arr[0,5] = 0.8853556084066838
$arr[0,6] = 0.9300466622665482
$ws.Range("C6:I6").Value = $arr
$ws.Range("N6").Value = 2.280343261403857

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.002555321430971702
$arr[0,1] = 0.02680069713119337
$arr[0,2] = 0.3425210282276652
$arr[0,3] = 1.253391139561273
$arr[0,4] = 1.28652295293557
$arr[0,5] = 0.905732788443288
$arr[0,6] = 0.9608679430435103
$ws.Range("C7:I7").Value = $arr
$ws.Range("N7").Value = 2.370273851392596

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.003104574511947789
$arr[0,1] = 0.03206250382704923
$arr[0,2] = 0.4138756201795815
$arr[0,3] = 1.430615228876633
$arr[0,4] = 1.501451066354917
$arr[0,5] = 0.9973129155281981
$arr[0,6] = 1.098511222843143
$ws.Range("C8:I8").Value = $arr
$ws.Range("N8").Value = 2.766433886209825

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.004193089834409136
$arr[0,1] = 0.04253537497490356
$arr[0,2] = 0.555641791949725
$arr[0,3] = 1.78643372609298
$arr[0,4] = 1.930770471933329
$arr[0,5] = 1.183194529484183
$arr[0,6] = 1.374773752377507
$ws.Range("C9:I9").Value = $arr
$ws.Range("N9").Value = 3.540180268007646

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.005002855156636343
$arr[0,1] = 0.05034663962427999
$arr[0,2] = 0.6612548373532121
$arr[0,3] = 2.053503134825348
$arr[0,4] = 2.251699954541436
$arr[0,5] = 1.323928293468327
$arr[0,6] = 1.582080857590313
$ws.Range("C10:I10").Value = $arr
$ws.Range("N10").Value = 4.10623028343673

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.00537410761167223
$arr[0,1] = 0.05393065389307594
$arr[0,2] = 0.7096928564080116
$arr[0,3] = 2.176361795961952
$arr[0,4] = 2.399058331659887
$arr[0,5] = 1.38893086000644
$arr[0,6] = 1.677437594168651
$ws.Range("C11:I11").Value = $arr
$ws.Range("N11").Value = 4.363110593465422

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.005515163530002809
$arr[0,1] = 0.05529264728627936
$arr[0,2] = 0.7280979333848023
$arr[0,3] = 2.22309264949709
$arr[0,4] = 2.455068630089386
$arr[0,5] = 1.413692841925524
$arr[0,6] = 1.713706417683738
$ws.Range("C12:I12").Value = $arr
$ws.Range("N12").Value = 4.460285735714251

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.005484762777086871
$arr[0,1] = 0.05499909744310116
$arr[0,2] = 0.7241311873572442
$arr[0,3] = 2.213018952958947
$arr[0,4] = 2.442996318239921
$arr[0,5] = 1.408353279212292
$arr[0,6] = 1.705888062043641
$ws.Range("C13:I13").Value = $arr
$ws.Range("N13").Value = 4.439361943450422

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.005385702608680987
$arr[0,1] = 0.05404260727628696
$arr[0,2] = 0.7112057626940071
$arr[0,3] = 2.180202163862333
$arr[0,4] = 2.403662073796795
$arr[0,5] = 1.390965065424609
$arr[0,6] = 1.680418212570942
$ws.Range("C14:I14").Value = $arr
$ws.Range("N14").Value = 4.371107314139522

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.005325088373172804
$arr[0,1] = 0.05345736725726624
$arr[0,2] = 0.7032969061931595
$arr[0,3] = 2.160128194102469
$arr[0,4] = 2.379596291095481
$arr[0,5] = 1.380333587760163
$arr[0,6] = 1.664838182651749
$ws.Range("C15:I15").Value = $arr
$ws.Range("N15").Value = 4.329286057409945

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.004978655728134385
$arr[0,1] = 0.05011306619284994
$arr[0,2] = 0.6580977365153871
$arr[0,3] = 2.045502423609747
$arr[0,4] = 2.242098281615483
$arr[0,5] = 1.319700468231588
$arr[0,6] = 1.575870917717282
$ws.Range("C16:I16").Value = $arr
$ws.Range("N16").Value = 4.089429168003562

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.004766908780190704
$arr[0,1] = 0.04806957234949039
$arr[0,2] = 0.6304746509789823
$arr[0,3] = 1.975540707466706
$arr[0,4] = 2.158106315897328
$arr[0,5] = 1.282759550197568
$arr[0,6] = 1.521567415056467
$ws.Range("C17:I17").Value = $arr
$ws.Range("N17").Value = 3.94211849063862

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.004645385761119059
$arr[0,1] = 0.04689706723252129
$arr[0,2] = 0.6146233109359684
$arr[0,3] = 1.935428655982889
$arr[0,4] = 2.109924162734671
$arr[0,5] = 1.261604197833208
$arr[0,6] = 1.490431971508826
$ws.Range("C18:I18").Value = $arr
$ws.Range("N18").Value = 3.857331695637754

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.004604284556602778
$arr[0,1] = 0.04650055608523473
$arr[0,2] = 0.6092624483873408
$arr[0,3] = 1.921869092488407
$arr[0,4] = 2.093632089949551
$arr[0,5] = 1.254457003917366
$arr[0,6] = 1.479906719779933
$ws.Range("C19:I19").Value = $arr
$ws.Range("N19").Value = 3.828614786364199

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.004789421446698583
$arr[0,1] = 0.04828680705101362
$arr[0,2] = 0.6334113316135728
$arr[0,3] = 1.982974922312167
$arr[0,4] = 2.167034087632373
$arr[0,5] = 1.286682402595602
$arr[0,6] = 1.527337864055909
$ws.Range("C20:I20").Value = $arr
$ws.Range("N20").Value = 3.95780600327754

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.005414785755526452
$arr[0,1] = 0.05432341797150286
$arr[0,2] = 0.7150005244286319
$arr[0,3] = 2.189835550321845
$arr[0,4] = 2.415209724981537
$arr[0,5] = 1.396068371685544
$arr[0,6] = 1.687894936110695
$ws.Range("C21:I21").Value = $arr
$ws.Range("N21").Value = 4.391158149571083

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.005826268009926139
$arr[0,1] = 0.0582968932927912
$arr[0,2] = 0.7686918527812026
$arr[0,3] = 2.326242151678684
$arr[0,4] = 2.578630384526264
$arr[0,5] = 1.468417554920961
$arr[0,6] = 1.793760671793933
$ws.Range("C22:I22").Value = $arr
$ws.Range("N22").Value = 4.67379181795809

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.005606379835651154
$arr[0,1] = 0.05617345905976379
$arr[0,2] = 0.7400000809948182
$arr[0,3] = 2.253325073928409
$arr[0,4] = 2.491293544283394
$arr[0,5] = 1.429722897787656
$arr[0,6] = 1.737170097874923
$ws.Range("C23:I23").Value = $arr
$ws.Range("N23").Value = 4.523002190001307

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.004779242810606377
$arr[0,1] = 0.04818858802545378
$arr[0,2] = 0.6320835667387712
$arr[0,3] = 1.979613572332795
$arr[0,4] = 2.162997512968161
$arr[0,5] = 1.284908624568857
$arr[0,6] = 1.524728781897011
$ws.Range("C24:I24").Value = $arr
$ws.Range("N24").Value = 3.950713976768498

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.003897082456433054
$arr[0,1] = 0.03968333960929726
$arr[0,2] = 0.5170591309376249
$arr[0,3] = 1.689227532280825
$arr[0,4] = 1.813714311031958
$arr[0,5] = 1.13220129104684
$arr[0,6] = 1.299310391377247
$ws.Range("C25:I25").Value = $arr
$ws.Range("N25").Value = 3.331249627311138
